$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (before) values for rows 2-4 in the affected columns
# (D, I, J, K, L, M, N, P, Q) before overwriting anything.
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

# Cyclic shift: new row2 = old row3, new row3 = old row4, new row4 = old row2
foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row3[$col]
    $ws.Range("${col}3").Value2 = $row4[$col]
    $ws.Range("${col}4").Value2 = $row2[$col]
}
